$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update symbol, fund name and classification for the sector that changed
# from Communication Services (XLC) to Energy (XLE)
$ws.Range("A5").Value = "XLE"
$ws.Range("B5").Value = "Energy Select Sector SPDR Fund"
$ws.Range("C5").Value = "Natural Resources Funds"

# Update disclaimer text with new "as of" date
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-06-09 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2500135127183522
$ws.Range("E2").Value = -0.007816091954023108

$ws.Range("D3").Value = 0.2499863450183538
$ws.Range("E3").Value = -0.009443861490031513

$ws.Range("D4").Value = 0.2499918152714616
$ws.Range("E4").Value = -0.01020213577421814

$ws.Range("D5").Value = 0.2500083269918323
$ws.Range("E5").Value = -0.006251116270762691

$ws.Range("D6").Value = 1
$ws.Range("E6").Value = -0.008428246584434973

# Re-apply sheet protection (best-effort; original used a legacy password hash
# that cannot be reproduced via the object model)
$ws.Protect()
